# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the "全部类型" sheet to reflect newer scrape data (gh-pages output
# regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 5570
$wsExpo.Range("F3").Value = 12881
$wsExpo.Range("F4").Value = 312
$wsExpo.Range("F5").Value = 629
$wsExpo.Range("F6").Value = 202
$wsExpo.Range("F7").Value = 392
$wsExpo.Range("F8").Value = 1207

# --- Sheet "全部类型" (all types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5570
$wsAll.Range("F4").Value = 12881
$wsAll.Range("F5").Value = 312
$wsAll.Range("F6").Value = 629
$wsAll.Range("F7").Value = 202
$wsAll.Range("F10").Value = 392
$wsAll.Range("F11").Value = 1207
